$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Assigning a literal string like "0%" through .Value triggers Excel's
    # smart "percentage" auto-detection, turning it into a numeric percent
    # cell with a new number format/style. Forcing the cell to Text format
    # first keeps it a literal string; resetting the style back to Normal
    # afterwards avoids leaving a stray style index behind.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Swap the header labels in D1 and E1
$ws.Range("D1").Value = "MORE THAN DEMAND PAYOUT%"
$ws.Range("E1").Value = "LESS THAN DEMAND PAYOUT%"

# Row 2
$ws.Range("C2").Value = 76668
Set-TextValue $ws.Range("D2") "0%"

# Row 3
$ws.Range("C3").Value = 109655
Set-TextValue $ws.Range("D3") "0%"
Set-TextValue $ws.Range("E3") "1%"

# Row 4
$ws.Range("C4").Value = 240187
Set-TextValue $ws.Range("D4") "0%"
Set-TextValue $ws.Range("E4") "2%"

# Row 5
$ws.Range("C5").Value = 80782

# Row 6
$ws.Range("C6").Value = 39800
Set-TextValue $ws.Range("E6") "0%"

# Row 7
$ws.Range("C7").Value = 149995
